# Update with Correct Forecast output
# - Rename Sheet1 -> "Sales vs PO"
# - Insert a new "Order Week" column (C) holding the original order dates,
#   shift "PO_Requested_Qty" to column D, and refresh the "ds" (A) column
#   with the (now later) week-ending dates.
# - Add three new sheets: "Weekly Growth", "Volume Insights", "Prediction Info"

$wb = $excel.ActiveWorkbook

# ---- Sheet1: rename + rebuild with the new "Order Week" column ----
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sales vs PO"

$ws1.Cells.Item(1,1).Value = "ds"
$ws1.Cells.Item(1,2).Value = "y"
$ws1.Cells.Item(1,3).Value = "Order Week"
$ws1.Cells.Item(1,4).Value = "PO_Requested_Qty"

# Make sure the newly-introduced D1 header cell (previously unused)
# picks up the same bold/centered/boxed header formatting as A1:C1.
$ws1.Range("A1").Copy()
$ws1.Range("D1").PasteSpecial(-4122)

# Each source row is (old ds date serial, y, PO_Requested_Qty).
# The new ds date is 6 days after the old date; the old date itself now
# lives in the new "Order Week" column.
$rows = @(
    @(44921,0,0),
    @(44935,0,0),
    @(44991,1,0),
    @(44998,1,0),
    @(45033,2,0),
    @(45040,0,0),
    @(45047,0,0),
    @(45054,0,0),
    @(45061,0,0),
    @(45068,1,0),
    @(45075,1,0),
    @(45299,0,0),
    @(45355,0,0),
    @(45362,0,0),
    @(45369,1,0),
    @(45376,1,0),
    @(45523,0,0),
    @(45558,0,0),
    @(45565,0,0),
    @(45572,2,0),
    @(45579,1,0),
    @(45586,1,0),
    @(45593,3,0),
    @(45600,0,0),
    @(45607,1,0),
    @(45614,0,0),
    @(45621,1,0)
)

$r = 2
foreach ($row in $rows) {
    $oldDs = $row[0]
    $y = $row[1]
    $poQty = $row[2]

    $ws1.Cells.Item($r,1).Value = $oldDs + 6
    $ws1.Cells.Item($r,1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws1.Cells.Item($r,2).Value = $y

    $ws1.Cells.Item($r,3).Value = $oldDs
    $ws1.Cells.Item($r,3).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws1.Cells.Item($r,4).Value = $poQty

    $r = $r + 1
}

# Reusable "header" format (bold, centered, boxed) taken from the
# existing Sheet1 header cells.
$headerFormat = $ws1.Range("A1")
$xlPasteFormats = -4122

# ---- New sheet: Weekly Growth ----
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Weekly Growth"
$ws2.Cells.Item(1,1).Value = "ds"
$ws2.Cells.Item(1,2).Value = "PO_Requested_Qty"
$ws2.Cells.Item(1,3).Value = "Growth%"
$headerFormat.Copy()
$ws2.Range("A1:C1").PasteSpecial($xlPasteFormats)

# ---- New sheet: Volume Insights ----
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "Volume Insights"
$ws3.Cells.Item(1,1).Value = "Total_PO_Quantity"
$ws3.Cells.Item(1,2).Value = "Average_PO_Quantity"
$ws3.Cells.Item(1,3).Value = "Max_PO_Quantity"
$ws3.Cells.Item(1,4).Value = "Min_PO_Quantity"
$ws3.Cells.Item(2,1).Value = 0
$ws3.Cells.Item(2,2).Value = 0
$ws3.Cells.Item(2,3).Value = 0
$ws3.Cells.Item(2,4).Value = 0
$headerFormat.Copy()
$ws3.Range("A1:D1").PasteSpecial($xlPasteFormats)

# ---- New sheet: Prediction Info ----
$ws4 = $wb.Worksheets.Add($null, $ws3)
$ws4.Name = "Prediction Info"
$ws4.Cells.Item(1,1).Value = "Predicted_Next_Week_PO_Quantity"
$ws4.Cells.Item(2,1).Value = 0
$headerFormat.Copy()
$ws4.Range("A1:A1").PasteSpecial($xlPasteFormats)

# Leave the first sheet active, matching the original workbook view.
$ws1.Activate()
